# MIL_2009.xlsx - add the team's 2009 season record (Wins / Losses / Ties)
# as three new trailing columns (AD, AE, AF) on the player-stats sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should look like the existing header row (bold, thin
# border, centered) - copy that formatting from the last existing header
# cell (AC1) onto AD1:AF1, then fill in the labels.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player (rows 2-47) shares the same team record: 80-82-0.
for ($row = 2; $row -le 47; $row++) {
    $ws.Range("AD$row").Value = 80
    $ws.Range("AE$row").Value = 82
    $ws.Range("AF$row").Value = 0
}
